$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "population census results" row (row 2) is no longer needed -
# remove it entirely, shifting the rest of the table up.
$ws.Rows("2:2").Delete()

# Only the 2014 figures are kept going forward; drop the 1989 and 2002
# columns (columns B and C), shifting the 2014 column into place as B.
$ws.Columns("B:C").Delete()

# Give the worksheet tab its proper municipality name instead of the
# generic "1".
$ws.Name = "გორი"

# Restore the active selection to A2, matching the saved view state.
[void]$ws.Range("A2").Select()
